# Updated cryptos list with GitHub Actions
# Apply updated Price (column D) and Volume(1h) (column E) values per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.657.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.808.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.806.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.452.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.819.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.689.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.738"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.94%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.135"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("E39").Value = "  -6.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.833.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0349"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.70%  "
